$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 150.16667
$ws.Range("I33").Value = 98.59999999999999
$ws.Range("J33").Value = 408
$ws.Range("K33").Value = 98.59999999999999
$ws.Range("L33").Value = 408
$ws.Range("M33").Value = 130.4
$ws.Range("N33").Value = -866
$ws.Range("H100").Value = 9261537
$ws.Range("I100").Value = 16668186
$ws.Range("J100").Value = 3225.75
$ws.Range("K100").Value = 16668186
$ws.Range("L100").Value = 3225.75
$ws.Range("M100").Value = -16667645
$ws.Range("N100").Value = -4307.75
$ws.Range("H106").Value = 15874960
$ws.Range("I106").Value = 15874960
$ws.Range("K106").Value = 15874960
$ws.Range("M106").Value = -15874329
$ws.Range("H138").Value = 11090795
$ws.Range("I138").Value = 1698562.9
$ws.Range("J138").Value = 17244326
$ws.Range("K138").Value = 5095688.699999999
$ws.Range("L138").Value = 51732978
$ws.Range("M138").Value = -5090548.699999999
$ws.Range("N138").Value = -51743258

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6055.8696
$ws.Range("I2").Value = 8226.375
$ws.Range("J2").Value = 1094.7142
$ws.Range("K2").Value = 8226.375
$ws.Range("L2").Value = 1094.7142
$ws.Range("M2").Value = -8113.375
$ws.Range("N2").Value = -1320.7142
$ws.Range("H32").Value = 3200.8604
$ws.Range("I32").Value = 2294.5151
$ws.Range("K32").Value = 2294.5151
$ws.Range("M32").Value = -2007.5151
$ws.Range("H74").Value = 11631.417
$ws.Range("I74").Value = 1807.7
$ws.Range("K74").Value = 1807.7
$ws.Range("M74").Value = -933.7
$ws.Range("H77").Value = 11631.417
$ws.Range("I77").Value = 1807.7
$ws.Range("K77").Value = 9038.5
$ws.Range("M77").Value = -4670.5
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1378
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 6055.8696
$ws.Range("I116").Value = 8226.375
$ws.Range("J116").Value = 1094.7142
$ws.Range("K116").Value = 8226.375
$ws.Range("L116").Value = 1094.7142
$ws.Range("M116").Value = -5932.375
$ws.Range("N116").Value = -5682.7142
$ws.Range("H123").Value = 33618.332
$ws.Range("J123").Value = 33618.332
$ws.Range("L123").Value = 33618.332
$ws.Range("N123").Value = -43418.332

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6055.8696
$ws.Range("I3").Value = 8226.375
$ws.Range("J3").Value = 1094.7142
$ws.Range("K3").Value = 8226.375
$ws.Range("L3").Value = 1094.7142
$ws.Range("M3").Value = -8112.375
$ws.Range("N3").Value = -1322.7142
$ws.Range("H94").Value = 885.5
$ws.Range("I94").Value = 866
$ws.Range("J94").Value = 1100
$ws.Range("K94").Value = 866
$ws.Range("L94").Value = 1100
$ws.Range("M94").Value = -415
$ws.Range("N94").Value = -2002
$ws.Range("H99").Value = 1228.125
$ws.Range("I99").Value = 1243.3334
$ws.Range("K99").Value = 1243.3334
$ws.Range("M99").Value = 254.6666
$ws.Range("H105").Value = 2950.5
$ws.Range("I105").Value = 2878.8262
$ws.Range("K105").Value = 2878.8262
$ws.Range("M105").Value = -1131.8262

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2952.9678
$ws.Range("I31").Value = 1251.925
$ws.Range("J31").Value = 6045.773
$ws.Range("K31").Value = 1251.925
$ws.Range("L31").Value = 6045.773
$ws.Range("M31").Value = -956.925
$ws.Range("N31").Value = -6635.773
$ws.Range("H34").Value = 2952.9678
$ws.Range("I34").Value = 1251.925
$ws.Range("J34").Value = 6045.773
$ws.Range("K34").Value = 1251.925
$ws.Range("L34").Value = 6045.773
$ws.Range("M34").Value = -1049.925
$ws.Range("N34").Value = -6449.773
$ws.Range("H105").Value = 528.1667
$ws.Range("I105").Value = 528.1667
$ws.Range("K105").Value = 528.1667
$ws.Range("M105").Value = 1218.8333

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2200
$ws.Range("I116").Value = 1250
$ws.Range("J116").Value = 2833.3333
$ws.Range("K116").Value = 3750
$ws.Range("L116").Value = 8499.999899999999
$ws.Range("M116").Value = -308
$ws.Range("N116").Value = -15383.9999
$ws.Range("H131").Value = 2599.9114
$ws.Range("J131").Value = 2599.9114
$ws.Range("L131").Value = 7799.7342
$ws.Range("N131").Value = -17879.7342
$ws.Range("H132").Value = 925.4666999999999
$ws.Range("I132").Value = 940.2857
$ws.Range("J132").Value = 912.5
$ws.Range("K132").Value = 8462.5713
$ws.Range("L132").Value = 8212.5
$ws.Range("M132").Value = -5932.5713
$ws.Range("N132").Value = -13272.5

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 3333.3333
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 6000
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -1673
$ws.Range("N55").Value = -6654
$ws.Range("H97").Value = 1087.7693
$ws.Range("I97").Value = 943
$ws.Range("J97").Value = 1570.3334
$ws.Range("K97").Value = 943
$ws.Range("L97").Value = 1570.3334
$ws.Range("M97").Value = -447
$ws.Range("N97").Value = -2562.3334
$ws.Range("H102").Value = 1732.8667
$ws.Range("I102").Value = 1290.4546
$ws.Range("J102").Value = 2949.5
$ws.Range("K102").Value = 1290.4546
$ws.Range("L102").Value = 2949.5
$ws.Range("M102").Value = 331.5454
$ws.Range("N102").Value = -6193.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2758.147
$ws.Range("I7").Value = 2043.8
$ws.Range("J7").Value = 3322.1052
$ws.Range("K7").Value = 2043.8
$ws.Range("L7").Value = 3322.1052
$ws.Range("M7").Value = -1931.8
$ws.Range("N7").Value = -3546.1052
$ws.Range("H100").Value = 2678.0344
$ws.Range("J100").Value = 3013.3333
$ws.Range("L100").Value = 3013.3333
$ws.Range("N100").Value = -4095.3333
$ws.Range("H126").Value = 2758.147
$ws.Range("I126").Value = 2043.8
$ws.Range("J126").Value = 3322.1052
$ws.Range("K126").Value = 6131.4
$ws.Range("L126").Value = 9966.3156
$ws.Range("M126").Value = -3661.4
$ws.Range("N126").Value = -14906.3156
$ws.Range("H131").Value = 38766.668
$ws.Range("I131").Value = 33150
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 33150
$ws.Range("L131").Value = 50000
$ws.Range("M131").Value = -28110
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 4572.421
$ws.Range("I132").Value = 3059.75
$ws.Range("J132").Value = 5672.5454
$ws.Range("K132").Value = 9179.25
$ws.Range("L132").Value = 17017.6362
$ws.Range("M132").Value = -6649.25
$ws.Range("N132").Value = -22077.6362

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 500
$ws.Range("I100").Value = 500
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -459
$ws.Range("N100").ClearContents()
$ws.Range("H123").Value = 38170.6
$ws.Range("J123").Value = 38170.6
$ws.Range("L123").Value = 38170.6
$ws.Range("N123").Value = -47970.6
$ws.Range("H126").Value = 72728.5
$ws.Range("I126").Value = 167650.17
$ws.Range("K126").Value = 502950.51
$ws.Range("M126").Value = -500480.51
$ws.Range("H132").Value = 20004940
$ws.Range("I132").Value = 33338950
$ws.Range("J132").Value = 3924.5
$ws.Range("K132").Value = 100016850
$ws.Range("L132").Value = 11773.5
$ws.Range("M132").Value = -100014320
$ws.Range("N132").Value = -16833.5
